# -----------------------------------------------------------------------
# Goal (per commit): add a new "Player Info" sheet (with player bio data)
# ahead of the existing "ODI Batting" sheet, and on "ODI Batting" turn the
# MATCH_CARD_LINK column (full scorecard URL) into a bare MATCH_CODE number.
#
# Note on sheet identity: in this COM shim, worksheet variables obtained via
# Worksheets.Item(<index>) are positional, not stable object handles - once
# a new sheet is inserted ahead of one, the old variable silently starts
# resolving to whatever now sits at that index. To dodge that entirely (and
# to land the new sheet on workbook sheetId=1 / the original data on
# sheetId=2, matching a natural "insert before" edit), we duplicate the
# existing sheet (the duplicate takes the new, higher sheetId) and turn the
# *duplicate* into the updated "ODI Batting", while repurposing the
# original in place into "Player Info". That also means the "Player Info"
# header row inherits the exact existing bold/border/center header style
# for free (no new style records), and the same trick (paste-special
# "Formats" from an always-plain cell) is used to strip the stray
# NumberFormat-only style that setting text-looking numbers leaves behind.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item(1)          # "ODI Batting", sheetId 1

# 1) Duplicate "ODI Batting" right after itself. The duplicate is a fresh
#    sheet object (new sheetId), carrying an exact copy of all data/styles.
$orig.Copy($null, $orig)
$dup = $wb.Worksheets.Item(2)
$dup.Name = "ODI Batting (tmp)"

# 2) Repurpose the duplicate into the updated "ODI Batting" sheet.
$dup.Range("D1").Value = "MATCH_CODE"

$dup.Range("D2").NumberFormat = "@"
$dup.Range("D2").Value = "4231"
$dup.Range("D3").NumberFormat = "@"
$dup.Range("D3").Value = "4232"
$dup.Range("D4").NumberFormat = "@"
$dup.Range("D4").Value = "4233"

# D2:D4 were plain/unstyled before; clear the NumberFormat-induced style we
# just introduced by cloning the (unstyled) format from A2 onto them.
$dup.Range("A2").Copy()
$dup.Range("D2:D4").PasteSpecial(-4122)

# 3) Repurpose the original sheet (sheetId 1) into "Player Info": drop the
#    now-unneeded extra columns/rows, keep the existing header style on
#    row 1, and lay down the new 4-column player-bio table.
$orig.Range("E1:J4").Clear()
$orig.Range("A3:D4").Clear()

$orig.Range("A1").Value = "ID"
$orig.Range("B1").Value = "NAME"
$orig.Range("C1").Value = "BATTING_HAND"
$orig.Range("D1").Value = "BOWL_STYLE"

$orig.Range("A2").NumberFormat = "@"
$orig.Range("A2").Value = "4732"
$orig.Range("B2").Value = "Tim Louis Seifert"
$orig.Range("C2").Value = "Right Handed"
$orig.Range("D2").Value = "Right Arm Fast Medium"

# A2 picked up a NumberFormat-only style from the text-forcing trick above;
# clear it the same way, borrowing the plain format already sitting on B2.
$orig.Range("B2").Copy()
$orig.Range("A2").PasteSpecial(-4122)

# 4) Names/order: original -> "Player Info" (before), duplicate -> "ODI
#    Batting" (after). Renaming the original only now avoids any name clash
#    with the duplicate's current temp name.
$orig.Name = "Player Info"
$dup.Name = "ODI Batting"

# Restore the original active sheet (first tab) so workbook-level view
# state (activeTab) is unchanged by these edits.
$orig.Select()
